$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Appends the new daily volatility-index readings (MOVE / VIX) that were
# published after "03-09-2021" (the previous last row, 172): five more
# trading dates running through "13-09-2021". The VIX figure for
# "13-09-2021" had already been published while MOVE had not yet been
# reported for that date, so column B is left blank on the final new row.

$newRows = @(
    @{ Date = "07-09-2021"; MOVE = "57.05"; VIX = "18.14" },
    @{ Date = "08-09-2021"; MOVE = "57.02"; VIX = "17.96" },
    @{ Date = "09-09-2021"; MOVE = "52.61"; VIX = "18.8" },
    @{ Date = "10-09-2021"; MOVE = "51.73"; VIX = "20.95" },
    @{ Date = "13-09-2021"; MOVE = $null;   VIX = "19.2" }
)

$startRow = 173

# Writing the date text directly through .Value would let the engine's
# smart "looks-like-a-date" detection silently coerce strings such as
# "07-09-2021" into a date serial (and stamp the cell with a new
# NumberFormat style). Routing the literal through a quoted formula and
# then collapsing it to a value with Copy + PasteSpecial(xlPasteValues)
# keeps the exact original text as a plain shared-string cell, with no
# style/number-format side effects - matching how the other ~170 date
# cells in this column are stored.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $entry = $newRows[$i]

    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.Formula = "=""" + $entry.Date + """"
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)

    if ($null -ne $entry.MOVE) {
        $ws.Cells.Item($row, 2).Value = [double]$entry.MOVE
    }

    $ws.Cells.Item($row, 3).Value = [double]$entry.VIX
}

$excel.CutCopyMode = $false
